$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132, shifting existing rows 132..231 down to 133..232.
$ws.Rows(132).Insert()

# Populate the newly inserted row 132 with the new weekly price record.
$ws.Range("A132").Value = 9
$ws.Range("B132").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C132").Value = "Metropolitana"
$ws.Range("D132").Value = 44582
$ws.Range("E132").Value = 13
$ws.Range("F132").Value = 100112001
$ws.Range("G132").Value = "Berenjena"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 52
$ws.Range("K132").Value = 10000
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = 11000
$ws.Range("N132").Value = "$/caja 60 unidades"
$ws.Range("O132").Value = "Región Metropolitana"
$ws.Range("P132").Value = 183
$ws.Range("Q132").Value = 60
$ws.Range("R132").Value = "Hortaliza"
